# Update "想去人数" (want-to-go count) figures across the workbook to the
# latest scrape values (gh-pages output regenerated at commit 456a3b4).
#
# The same events appear on multiple sheets:
#   展览      (sheet1, "Exhibitions")
#   演出      (sheet2, "Performances")
#   本地生活  (sheet3, "Local life")
#   全部类型  (sheet4, "All types" — union of the above three)
# so each updated count is written to every sheet/row where that event shows
# up, keeping column F in sync everywhere.

$wb = $excel.ActiveWorkbook

function Set-F {
    param($SheetName, $Row, $NewValue)
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Cells.Item($Row, 6).Value = $NewValue
}

# --- 展览 ---
Set-F "展览" 3  1129
Set-F "展览" 6  467
Set-F "展览" 7  755
Set-F "展览" 10 48
Set-F "展览" 11 420
Set-F "展览" 14 893
Set-F "展览" 16 2021
Set-F "展览" 17 511
Set-F "展览" 18 8020
Set-F "展览" 19 609
Set-F "展览" 21 67
Set-F "展览" 22 99

# --- 演出 ---
Set-F "演出" 2 525
Set-F "演出" 6 12

# --- 本地生活 ---
Set-F "本地生活" 2 5568
Set-F "本地生活" 4 396

# --- 全部类型 ---
Set-F "全部类型" 3  5568
Set-F "全部类型" 5  396
Set-F "全部类型" 6  525
Set-F "全部类型" 7  1129
Set-F "全部类型" 12 467
Set-F "全部类型" 13 755
Set-F "全部类型" 17 12
Set-F "全部类型" 18 48
Set-F "全部类型" 19 420
Set-F "全部类型" 24 893
Set-F "全部类型" 28 2021
Set-F "全部类型" 29 511
Set-F "全部类型" 30 8020
Set-F "全部类型" 33 609
Set-F "全部类型" 35 67
Set-F "全部类型" 36 99
